$d = $word.ActiveDocument

# Fill in the "Name:" placeholder on the cover section with the student's
# actual name, and drop the yellow highlight that was marking it as a
# to-do placeholder.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "[Insert your name here]"
$find.Forward = $true
$find.Wrap = 1

$found = $find.Execute()
if ($found) {
    $rng = $find.Parent
    $rng.Text = "Tracy Robert Mann"
    $rng.HighlightColorIndex = 0
}
